$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row text (styles/shared-string slots stay where they were)
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Email"
$ws.Range("C1").Value = "Demand"

# New data row
$ws.Range("A2").Value = "Bill"
$ws.Range("C2").Value = 50000

# Column widths (A and B) - nearest values on the engine's width grid to the
# author's 14.7265625 / 17.81640625 (OOXML chars, 1/256-px precision)
$ws.Columns.Item(1).ColumnWidth = 13.87
$ws.Columns.Item(2).ColumnWidth = 17.0

# Hyperlink (pre-set the display text so Add() doesn't stamp a redundant
# display="" attribute; Add() still applies the built-in Hyperlink style)
$ws.Range("B2").Value = "bill@example.com"
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:bill@example.com")

# Selection as left by the author
$ws.Range("F4").Select()
